$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New list of URLs for rows 2..34 (A1 keeps header "link")
$urls = @(
    "https://www.facebook.com/K14vn/posts/967228682104798?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967221998772133?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967207705440229?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967186722108994?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967169058777427?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967148892112777?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967114128782920?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967100842117582?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967043255456674?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967033045457695?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/967009495460050?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966942885466711?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966706425490357?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966653402162326?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966566715504328?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966521775508822?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966476808846652?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966388225522177?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966375708856762?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966335782194088?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966306918863641?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966263705534629?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/965422742285392?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966245622203104?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966243385536661?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966202212207445?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966189408875392?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966167398877593?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/966134648880868?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/965910925569907?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/965679468926386?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/965634952264171?ref=embed_post",
    "https://www.facebook.com/K14vn/posts/965586232269043?ref=embed_post"
)

# Remove the old rows below the new data range (34 rows incl. header) first,
# so the used range/dimension shrinks to A1:A34.
$ws.Rows("35:84").Delete()

# Overwrite A2:A34 with the refreshed post links.
for ($i = 0; $i -lt $urls.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $urls[$i]
}

# Drop the stray per-cell formatting on A17 (it no longer carries a distinct style).
$ws.Range("A17").ClearFormats()

# Reset the view: no frozen/scrolled top-left row, selection parked on the last cell.
[void]$ws.Range("A34").Select()
